$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.955.91"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "3.381.84"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.06"
$ws.Range("E5").Value = "  +3.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.52"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("D8").Value = "3.371.61"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +12.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.631"
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.16"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("E13").Value = "  +6.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("E14").Value = "  +3.84%  "
$ws.Range("D15").Value = "3.919.08"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.25"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "3.389.07"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").Value = "64.929.26"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.75"
$ws.Range("E22").Value = "  +14.78%  "
$ws.Range("E23").Value = "  +14.02%  "
$ws.Range("E24").Value = "  +3.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.80"
$ws.Range("E25").Value = "  +5.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.53"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +7.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.79"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.77"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.74"
$ws.Range("E30").Value = "  +7.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  +6.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.51"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "569.52"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  +6.63%  "
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.67"
$ws.Range("E37").Value = "  +8.60%  "
$ws.Range("E38").Value = "  -3.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.60"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").Value = "3.086.46"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.86"
$ws.Range("E44").Value = "  +4.49%  "
$ws.Range("E45").Value = "  +4.83%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.134"
$ws.Range("E47").Value = "  +6.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.13"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.20"
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.28"
$ws.Range("E51").Value = "  +4.24%  "
